$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 650
$ws.Range("B3").Value = 343
$ws.Range("B4").Value = 78
$ws.Range("B5").Value = 64
$ws.Range("B6").Value = 300
